$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$updates = @(
    @(2, 6, 2637),
    @(4, 6, 468),
    @(6, 6, 198),
    @(7, 6, 480),
    @(8, 6, 1219),
    @(9, 6, 563),
    @(10, 6, 306),
    @(11, 3, "杭州·萌忧 动漫游戏嘉年华"),
    @(11, 4, "康候圣街99号 顺丰创新中心"),
    @(11, 5, "2024.07.06 10:30-07.06 17:00"),
    @(11, 6, 0),
    @(11, 7, 55),
    @(11, 8, "https://show.bilibili.com/platform/detail.html?id=87292"),
    @(11, 9, "//i0.hdslb.com/bfs/openplatform/202406/s4DzfyI11718083528496.jpeg"),
    @(12, 3, "杭州·重逢·怀旧only"),
    @(12, 4, "丰庆路492号建冠龙禾商务中心A幢 杭州华礼宴国际礼宴中心(龙禾商务中心店)"),
    @(12, 5, "2024.07.06 09:00-07.06 17:00"),
    @(12, 6, 123),
    @(12, 7, 69),
    @(12, 8, "https://show.bilibili.com/platform/detail.html?id=85742"),
    @(12, 9, "//i2.hdslb.com/bfs/openplatform/202405/qBeP0pEz1715399357252.jpeg"),
    @(13, 2, "2024-07-06"),
    @(13, 3, "杭州·黑执事only"),
    @(13, 4, "大岭山路156号 爱丽芬城堡"),
    @(13, 5, "2024.07.06 10:00-07.07 18:00"),
    @(13, 6, 357),
    @(13, 7, 160),
    @(13, 8, "https://show.bilibili.com/platform/detail.html?id=86414"),
    @(13, 9, "//i1.hdslb.com/bfs/openplatform/202405/iP2cxk2w1716800288950.jpeg"),
    @(14, 3, "杭州·AD04动漫展"),
    @(14, 5, "2024.07.13 10:00-07.14 17:00"),
    @(14, 6, 5703),
    @(14, 7, 75),
    @(14, 8, "https://show.bilibili.com/platform/detail.html?id=85012"),
    @(14, 9, "//i0.hdslb.com/bfs/openplatform/202405/y1iKqqnh1715326769523.jpeg"),
    @(15, 6, 1765),
    @(16, 6, 4132),
    @(17, 6, 428),
    @(19, 6, 303),
    @(20, 6, 4822),
    @(21, 6, 6214),
    @(23, 6, 1055),
    @(24, 6, 688),
    @(25, 6, 3763),
    @(26, 6, 496),
    @(27, 6, 67),
    @(28, 6, 190),
    @(29, 6, 129),
    @(30, 6, 985),
    @(31, 6, 1411),
    @(32, 6, 466),
    @(33, 6, 549),
    @(34, 6, 1595),
    @(35, 6, 200),
    @(36, 6, 1714),
    @(37, 6, 190),
    @(39, 6, 1133),
    @(40, 3, "杭州·梦漫星河动漫嘉年华·赵路专场"),
    @(40, 4, "阳城路雅澳杭州电商产业园西侧约200米 杭州大会展中心"),
    @(40, 5, "2024.08.04 11:40-08.04 17:00"),
    @(40, 6, 1336),
    @(40, 7, "已售罄"),
    @(40, 8, "https://show.bilibili.com/platform/detail.html?id=86221"),
    @(40, 9, "//i1.hdslb.com/bfs/openplatform/202405/2padflbr1716372780297.jpeg"),
    @(41, 2, "2024-08-10"),
    @(41, 3, "杭州·原神X星铁X绝区零only"),
    @(41, 4, "望江东路333号 杭州瑞莱克斯大酒店"),
    @(41, 5, "2024.08.10 10:00-08.10 17:00"),
    @(41, 6, 629),
    @(41, 7, 60),
    @(41, 8, "https://show.bilibili.com/platform/detail.html?id=82754"),
    @(41, 9, "//i1.hdslb.com/bfs/openplatform/202403/qA0LNJuF1710234461030.jpeg"),
    @(42, 3, "杭州·造梦探险家城堡二次元同好会"),
    @(42, 4, "大岭山路156号 爱丽芬城堡"),
    @(42, 5, "2024.08.10 10:00-08.10 22:00"),
    @(42, 6, 95),
    @(42, 7, 38),
    @(42, 8, "https://show.bilibili.com/platform/detail.html?id=86432"),
    @(42, 9, "//i2.hdslb.com/bfs/openplatform/202405/xWUy30Ns1716783723057.jpeg"),
    @(43, 2, "2024-08-17"),
    @(43, 3, "杭州·HD·01"),
    @(43, 4, "钱江世纪城奔竞大道353号 杭州国际博览中心"),
    @(43, 5, "2024.08.17 09:30-08.18 17:00"),
    @(43, 6, 3407),
    @(43, 7, 75),
    @(43, 8, "https://show.bilibili.com/platform/detail.html?id=86332"),
    @(43, 9, "//i2.hdslb.com/bfs/openplatform/202405/GBMur4hT1716145118862.jpeg"),
    @(44, 3, "浙江·马娘ONLY03-晴风游憩"),
    @(44, 4, "康候圣街99号 顺丰创新中心"),
    @(44, 5, "2024.08.17 10:00-08.17 17:00"),
    @(44, 6, 130),
    @(44, 7, 68),
    @(44, 8, "https://show.bilibili.com/platform/detail.html?id=86529"),
    @(44, 9, "//i1.hdslb.com/bfs/openplatform/202405/21d6moub1716799089058.jpeg"),
    @(45, 2, "2024-08-18"),
    @(45, 3, "浙江·蔚蓝档案ONLY02-夏末狂欢！"),
    @(45, 5, "2024.08.18 10:00-08.18 17:00"),
    @(45, 6, 284),
    @(45, 8, "https://show.bilibili.com/platform/detail.html?id=86594"),
    @(45, 9, "//i1.hdslb.com/bfs/openplatform/202405/TVvJFURG1716799911888.jpeg"),
    @(46, 2, "2024-08-24"),
    @(46, 3, "杭州·D3动漫游戏嘉年华"),
    @(46, 4, "德胜东路2539号 梦马汽车小镇"),
    @(46, 5, "2024.08.24 10:00-08.24 17:00"),
    @(46, 6, 409),
    @(46, 7, 50),
    @(46, 8, "https://show.bilibili.com/platform/detail.html?id=84912"),
    @(46, 9, "//i0.hdslb.com/bfs/openplatform/202405/tAGUPfkr1715667000279.jpeg"),
    @(47, 3, "杭州·萌忧·原崩铁only"),
    @(47, 4, "康候圣街99号 顺丰创新中心"),
    @(47, 5, "2024.08.24 10:30-08.24 17:00"),
    @(47, 6, 1),
    @(47, 8, "https://show.bilibili.com/platform/detail.html?id=87293"),
    @(47, 9, "//i2.hdslb.com/bfs/openplatform/202406/rQFz5smR1717475284585.jpeg"),
    @(48, 6, 11),
    @(49, 6, 3883)
)
foreach ($u in $updates) {
    $ws.Cells.Item($u[0], $u[1]).Value = $u[2]
}

$ws = $wb.Worksheets.Item("全部类型")
$updates = @(
    @(2, 6, 3897),
    @(4, 6, 2637),
    @(6, 6, 468),
    @(8, 6, 1201),
    @(10, 6, 198),
    @(11, 6, 480),
    @(13, 6, 1219),
    @(14, 6, 564),
    @(15, 6, 306),
    @(19, 6, 1765),
    @(20, 6, 4822),
    @(22, 6, 1055),
    @(23, 6, 688),
    @(24, 6, 3763),
    @(25, 6, 496),
    @(26, 6, 67),
    @(27, 6, 190),
    @(28, 6, 129),
    @(29, 6, 985),
    @(30, 6, 1411),
    @(31, 6, 466),
    @(32, 6, 550),
    @(34, 6, 1595),
    @(35, 6, 200),
    @(36, 6, 1714),
    @(39, 6, 630),
    @(41, 6, 95),
    @(42, 6, 74),
    @(43, 6, 3407),
    @(45, 6, 130),
    @(49, 6, 3883)
)
foreach ($u in $updates) {
    $ws.Cells.Item($u[0], $u[1]).Value = $u[2]
}

$ws = $wb.Worksheets.Item("演出")
$updates = @(
    @(5, 6, 1201),
    @(7, 6, 41),
    @(25, 6, 74)
)
foreach ($u in $updates) {
    $ws.Cells.Item($u[0], $u[1]).Value = $u[2]
}

$ws = $wb.Worksheets.Item("本地生活")
$updates = @(
    @(2, 6, 3897)
)
foreach ($u in $updates) {
    $ws.Cells.Item($u[0], $u[1]).Value = $u[2]
}
